# Weekly price-sheet refresh for "Hortaliza, Comercializadora del Agro de
# Limari - Arveja Verde": a new week of observations is inserted, which
# (per the source process) shifts every existing data row down by one
# position and appends the two rows that were pushed past the prior end
# of the table (old rows 43 and 44) as new rows 45-46.
#
# Net effect vs. the starting workbook, expressed as a direct cell-value
# patch (equivalent to, but simpler/more robust than, scripting 40+
# individual row-insert operations via COM):
#   - Rows 5-44: columns D (Fecha), J (Volumen), K/L/M (Precio min/max/prom)
#     and P (Precio $/Kg) take on the values the diff specifies.
#   - Rows 45-46 are brand-new rows, fully populated, carrying the same
#     "Comercializadora del Agro de Limari" / Arveja Verde schema as the
#     rest of the table.
#   - The sheet dimension grows from A1:R44 to A1:R46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows 5-44: update the columns the diff touches (Fecha/Volumen/Precios)
# ---------------------------------------------------------------------
# Row 5
$ws.Range("D5").Value = 44503
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 11000
$ws.Range("L5").Value = 13000
$ws.Range("M5").Value = 12000
$ws.Range("P5").Value = 480
# Row 6
$ws.Range("D6").Value = 44406
$ws.Range("K6").Value = 26000
$ws.Range("L6").Value = 28000
$ws.Range("M6").Value = 27000
$ws.Range("P6").Value = 1080
# Row 7
$ws.Range("D7").Value = 44419
$ws.Range("J7").Value = 600
$ws.Range("K7").Value = 27000
$ws.Range("L7").Value = 29000
$ws.Range("M7").Value = 28000
$ws.Range("P7").Value = 1120
# Row 8
$ws.Range("D8").Value = 44475
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = 22000
$ws.Range("L8").Value = 24000
$ws.Range("M8").Value = 23000
$ws.Range("P8").Value = 920
# Row 9
$ws.Range("D9").Value = 44363
$ws.Range("J9").Value = 240
$ws.Range("K9").Value = 28000
$ws.Range("L9").Value = 30000
$ws.Range("M9").Value = 29000
$ws.Range("P9").Value = 1160
# Row 10
$ws.Range("D10").Value = 44391
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 26000
$ws.Range("L10").Value = 28000
$ws.Range("M10").Value = 27000
$ws.Range("P10").Value = 1080
# Row 11
$ws.Range("D11").Value = 44364
$ws.Range("J11").Value = 200
# Row 12
$ws.Range("D12").Value = 44434
$ws.Range("J12").Value = 500
# Row 13
$ws.Range("D13").Value = 44441
$ws.Range("K13").Value = 28000
$ws.Range("L13").Value = 30000
$ws.Range("M13").Value = 29000
$ws.Range("P13").Value = 1160
# Row 14
$ws.Range("D14").Value = 44413
$ws.Range("J14").Value = 700
$ws.Range("K14").Value = 26000
$ws.Range("L14").Value = 28000
$ws.Range("M14").Value = 27000
$ws.Range("P14").Value = 1080
# Row 15
$ws.Range("D15").Value = 44448
# Row 16
$ws.Range("D16").Value = 44426
$ws.Range("K16").Value = 28000
$ws.Range("L16").Value = 30000
$ws.Range("M16").Value = 29000
$ws.Range("P16").Value = 1160
# Row 17
$ws.Range("D17").Value = 44489
$ws.Range("J17").Value = 400
$ws.Range("K17").Value = 18000
$ws.Range("L17").Value = 20000
$ws.Range("M17").Value = 19000
$ws.Range("P17").Value = 760
# Row 18
$ws.Range("D18").Value = 44461
$ws.Range("J18").Value = 500
$ws.Range("K18").Value = 23000
$ws.Range("L18").Value = 25000
$ws.Range("M18").Value = 24000
$ws.Range("P18").Value = 960
# Row 19
$ws.Range("D19").Value = 44455
$ws.Range("J19").Value = 800
# Row 20
$ws.Range("D20").Value = 44447
$ws.Range("J20").Value = 600
$ws.Range("K20").Value = 28000
$ws.Range("L20").Value = 30000
$ws.Range("M20").Value = 29000
$ws.Range("P20").Value = 1160
# Row 21
$ws.Range("D21").Value = 44483
$ws.Range("J21").Value = 300
$ws.Range("K21").Value = 18000
$ws.Range("L21").Value = 20000
$ws.Range("M21").Value = 19000
$ws.Range("P21").Value = 760
# Row 22
$ws.Range("D22").Value = 44435
$ws.Range("J22").Value = 900
# Row 23
$ws.Range("D23").Value = 44427
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 28000
$ws.Range("L23").Value = 30000
$ws.Range("M23").Value = 29000
$ws.Range("P23").Value = 1160
# Row 24
$ws.Range("D24").Value = 44468
$ws.Range("J24").Value = 500
$ws.Range("K24").Value = 23000
$ws.Range("L24").Value = 25000
$ws.Range("M24").Value = 24000
$ws.Range("P24").Value = 960
# Row 25
$ws.Range("D25").Value = 44420
$ws.Range("J25").Value = 700
$ws.Range("K25").Value = 27000
$ws.Range("L25").Value = 29000
$ws.Range("M25").Value = 28000
$ws.Range("P25").Value = 1120
# Row 26
$ws.Range("D26").Value = 44412
$ws.Range("J26").Value = 600
$ws.Range("K26").Value = 25000
$ws.Range("L26").Value = 27000
$ws.Range("M26").Value = 26000
$ws.Range("P26").Value = 1040
# Row 27
$ws.Range("D27").Value = 44377
# Row 28
$ws.Range("D28").Value = 44405
$ws.Range("J28").Value = 500
$ws.Range("K28").Value = 26000
$ws.Range("L28").Value = 28000
$ws.Range("M28").Value = 27000
$ws.Range("P28").Value = 1080
# Row 29
$ws.Range("D29").Value = 44350
$ws.Range("J29").Value = 700
$ws.Range("K29").Value = 28000
$ws.Range("L29").Value = 30000
$ws.Range("M29").Value = 29000
$ws.Range("P29").Value = 1160
# Row 30
$ws.Range("D30").Value = 44385
# Row 31
$ws.Range("D31").Value = 44398
$ws.Range("K31").Value = 26000
$ws.Range("L31").Value = 28000
$ws.Range("M31").Value = 27000
$ws.Range("P31").Value = 1080
# Row 32
$ws.Range("D32").Value = 44371
$ws.Range("J32").Value = 500
# Row 33
$ws.Range("D33").Value = 44454
$ws.Range("J33").Value = 1000
$ws.Range("K33").Value = 28000
$ws.Range("L33").Value = 30000
$ws.Range("M33").Value = 29000
$ws.Range("P33").Value = 1160
# Row 34
$ws.Range("D34").Value = 44490
$ws.Range("J34").Value = 500
$ws.Range("K34").Value = 16000
$ws.Range("L34").Value = 18000
$ws.Range("M34").Value = 17000
$ws.Range("P34").Value = 680
# Row 35
$ws.Range("D35").Value = 44399
$ws.Range("J35").Value = 400
$ws.Range("K35").Value = 26000
$ws.Range("L35").Value = 28000
$ws.Range("M35").Value = 27000
$ws.Range("P35").Value = 1080
# Row 36
$ws.Range("D36").Value = 44357
$ws.Range("J36").Value = 340
$ws.Range("K36").Value = 28000
$ws.Range("L36").Value = 30000
$ws.Range("M36").Value = 29000
$ws.Range("P36").Value = 1160
# Row 37
$ws.Range("D37").Value = 44476
$ws.Range("K37").Value = 23000
$ws.Range("L37").Value = 24000
$ws.Range("M37").Value = 23500
$ws.Range("P37").Value = 940
# Row 38
$ws.Range("D38").Value = 44482
$ws.Range("J38").Value = 500
$ws.Range("K38").Value = 18000
$ws.Range("L38").Value = 20000
$ws.Range("M38").Value = 19000
$ws.Range("P38").Value = 760
# Row 39
$ws.Range("D39").Value = 44356
$ws.Range("J39").Value = 300
$ws.Range("K39").Value = 26000
$ws.Range("L39").Value = 28000
$ws.Range("M39").Value = 27000
$ws.Range("P39").Value = 1080
# Row 40
$ws.Range("D40").Value = 44469
$ws.Range("J40").Value = 600
$ws.Range("K40").Value = 22000
$ws.Range("L40").Value = 24000
$ws.Range("M40").Value = 23000
$ws.Range("P40").Value = 920
# Row 41
$ws.Range("D41").Value = 44504
$ws.Range("J41").Value = 600
$ws.Range("K41").Value = 11000
$ws.Range("L41").Value = 13000
$ws.Range("M41").Value = 12000
$ws.Range("P41").Value = 480
# Row 42
$ws.Range("D42").Value = 44384
$ws.Range("K42").Value = 26000
$ws.Range("M42").Value = 27000
$ws.Range("P42").Value = 1080
# Row 43
$ws.Range("D43").Value = 44343
$ws.Range("J43").Value = 200
# Row 44
$ws.Range("D44").Value = 44370
$ws.Range("K44").Value = 27000
$ws.Range("L44").Value = 28000
$ws.Range("M44").Value = 27500
$ws.Range("P44").Value = 1100

# ---------------------------------------------------------------------
# Rows 45-46: new rows appended at the bottom of the table
# ---------------------------------------------------------------------
# Row 45
$ws.Range("A45").Value = 2
$ws.Range("B45").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C45").Value = "Coquimbo"
$ws.Range("D45").Value = 44392
$ws.Range("E45").Value = 4
$ws.Range("F45").Value = 100112022
$ws.Range("G45").Value = "Arveja Verde"
$ws.Range("H45").Value = "Perfection"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 100
$ws.Range("K45").Value = 26000
$ws.Range("L45").Value = 28000
$ws.Range("M45").Value = 27000
$ws.Range("N45").Value = "`$/malla 25 kilos"
$ws.Range("O45").Value = "Provincia de Limarí"
$ws.Range("P45").Value = 1080
$ws.Range("Q45").Value = 25
$ws.Range("R45").Value = "Hortaliza"
$ws.Range("D45").NumberFormat = $ws.Range("D4").NumberFormat()

# Row 46
$ws.Range("A46").Value = 2
$ws.Range("B46").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C46").Value = "Coquimbo"
$ws.Range("D46").Value = 44433
$ws.Range("E46").Value = 4
$ws.Range("F46").Value = 100112022
$ws.Range("G46").Value = "Arveja Verde"
$ws.Range("H46").Value = "Perfection"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 400
$ws.Range("K46").Value = 28000
$ws.Range("L46").Value = 30000
$ws.Range("M46").Value = 29000
$ws.Range("N46").Value = "`$/malla 25 kilos"
$ws.Range("O46").Value = "Provincia de Limarí"
$ws.Range("P46").Value = 1160
$ws.Range("Q46").Value = 25
$ws.Range("R46").Value = "Hortaliza"
$ws.Range("D46").NumberFormat = $ws.Range("D4").NumberFormat()
